$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "Dacă poți, ajută-i cu rutina zilnică, de exemplu, luați masa împreună și să faceți curățenie. ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Dacă poți, ajută-i cu rutina zilnică, de exemplu, luați masa împreună și faceți curățenie. ",
    2)

$d.Content.Find.Execute(
    "Încearcă să îți rezervi timp în fiecare zi pentru a te juca sau a petrece timp de calitate alături copiii. ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Încearcă să îți rezervi timp în fiecare zi pentru a te juca sau a petrece timp de calitate alături copii. ",
    2)

$d.Content.Find.Execute(
    "Aprecierile le arată copiilor că îi observi și că îți pasă.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Aprecierile le arată copiilor că îi vezi și că îți pasă.",
    2)

$d.Content.Find.Execute(
    "Să te gândești și măcar la un singur lucru de care ești mândru poate face diferența!",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Să te gândești măcar și la un singur lucru de care ești mândru poate face diferența!",
    2)
